$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-8: Matricula_Historico (B), Curso (D), Status (E), Voto Valido (G) ---

$matriculas = @("2022004547", "2023001138", "2019018872", "2023000239", "2022004420", "2023010450", "2023004087")
for ($i = 0; $i -lt $matriculas.Length; $i++) {
    $r = $i + 2
    $cellB = $ws.Cells.Item($r, 2)
    # Write the digit-string via a TEXT() formula, then freeze it to a static
    # value with PasteSpecial so the cell ends up holding plain text without
    # Excel's auto number-detection converting it to a numeric cell, and
    # without requiring a NumberFormat change (keeps the original style).
    $cellB.Formula = '=TEXT(' + $matriculas[$i] + ',"0")'
    $cellB.Copy()
    $cellB.PasteSpecial(-4163)
    $ws.Cells.Item($r, 4).Value = "SISTEMAS DE INFORMAÇÃO/CAMP/CAMB"
    $ws.Cells.Item($r, 5).Value = "ATIVO"
    $ws.Cells.Item($r, 7).Value = "Voto Valido."
}

# --- Row 9: only the Curso (D9) column changes ---
$ws.Cells.Item(9, 4).Value = "Arquivo Invalido"

# --- New row 10: separator row ---
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "-"
$ws.Cells.Item(10, 3).Value = "-"
$ws.Cells.Item(10, 4).Value = "-"
$ws.Cells.Item(10, 5).Value = "-"
$ws.Cells.Item(10, 6).Value = "-"
$ws.Cells.Item(10, 7).Value = "-"

# --- New row 11: summary row ---
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "Total de Votos Validos"
$ws.Cells.Item(11, 3).Value = 7
$ws.Cells.Item(11, 4).Value = "Votos por Chapa"
$ws.Cells.Item(11, 5).Value = "{'Nova Era': 6, 'Branco': 1}"
$ws.Cells.Item(11, 6).Value = "Nova Era"
$ws.Cells.Item(11, 7).Value = "Chapa Vencedora"

# --- Copy the style from A9 (existing bordered/bold/centered style) to A10 and A11 ---
$ws.Range("A9").Copy()
$ws.Range("A10:A11").PasteSpecial(-4122)
